# Update countries & provincias Spain
# Refresh the COVID-19 dashboard data to the 29 Abril 2020 19:52 snapshot.
#
# The source feed re-ranks countries by "Casos totales" (column B) in
# descending order, so as raw counts change a few countries swap places
# in the table: Pakistan overtakes Austria (rows 29/30) and Irak overtakes
# Uzbekistan / Afganistan / Armenia (rows 68-71). Everywhere else only the
# numeric columns for the affected rows are refreshed in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp banner (A1) -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 19:52"

# --- Row 4: Estados Unidos (simple in-place refresh) ------------------
$ws.Range("B4").Value = 1048834
$ws.Range("C4").Value = 13069
$ws.Range("E4").Value = 843987
$ws.Range("G4").Value = 1229
$ws.Range("H4").Value = 60495

# --- Row 14: Brasil (simple in-place refresh) --------------------------
$ws.Range("B14").Value = 74493
$ws.Range("C14").Value = 1594
$ws.Range("E14").Value = 36791
$ws.Range("G14").Value = 95
$ws.Range("H14").Value = 5158

# --- Rows 29-30: Pakistan overtakes Austria -----------------------------
$ws.Range("A29").Value = "Pakistan"
$ws.Range("B29").Value = 15525
$ws.Range("C29").Value = 913
$ws.Range("D29").Value = 3425
$ws.Range("E29").Value = 11757
$ws.Range("F29").Value = 111
$ws.Range("G29").Value = 31
$ws.Range("H29").Value = 343

$ws.Range("A30").Value = "Austria"
$ws.Range("B30").Value = 15402
$ws.Range("C30").Value = 45
$ws.Range("D30").Value = 12779
$ws.Range("E30").Value = 2043
$ws.Range("F30").Value = 131
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = 580

# --- Row 63: Barein (simple in-place refresh) ---------------------------
$ws.Range("B63").Value = 2869
$ws.Range("C63").Value = 58
$ws.Range("D63").Value = 1370
$ws.Range("E63").Value = 1491

# --- Rows 68-71: Irak overtakes Uzbekistan / Afganistan / Armenia ------
$ws.Range("A68").Value = "Irak"
$ws.Range("B68").Value = 2003
$ws.Range("C68").Value = 75
$ws.Range("D68").Value = 1346
$ws.Range("E68").Value = 565
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 92

$ws.Range("A69").Value = "Uzbekistan"
$ws.Range("B69").Value = 1969
$ws.Range("C69").Value = 30
$ws.Range("D69").Value = 1096
$ws.Range("E69").Value = 865
$ws.Range("F69").Value = 8
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 8

$ws.Range("A70").Value = "Afganistan"
$ws.Range("B70").Value = 1939
$ws.Range("C70").Value = 111
$ws.Range("D70").Value = 252
$ws.Range("E70").Value = 1627
$ws.Range("F70").Value = 7
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 60

$ws.Range("A71").Value = "Armenia"
$ws.Range("B71").Value = 1932
$ws.Range("C71").Value = 65
$ws.Range("D71").Value = 900
$ws.Range("E71").Value = 1002
$ws.Range("F71").Value = 10
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 30

# --- Row 152: Monaco (simple in-place refresh) --------------------------
$ws.Range("D152").Value = 58
$ws.Range("E152").Value = 33

Write-Output "Countries & provincias Spain data refreshed"
